$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the new phone number value into A7 (Excel will add it to the shared strings table)
$ws.Range("A7").Value = "(726) 222 - 1745"

# Update the selection to match the target state
$ws.Range("D9").Select()

# Column A was "best fit" before; re-fit it now that longer content was added
$ws.Columns("A:A").AutoFit()
